$wb = $excel.ActiveWorkbook

# Values to write into C11:I11 (GDP, UEMP, CPI, LTRate, EURUSD, WTI, RPP)
$cols = @("C", "D", "E", "F", "G", "H", "I")
$values = @(0.861662757638527, -0.8902777777777704, 0.29771488471092766, -0.14700000000000002, 1.5829618029997903, 16.12947350163202, 0.383631713554976)

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + "11").Value = $values[$i]
    }
}
